$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

$old = "go to week 5"
$new = "download week 5"
$fullText = $tr.Text
$startPos = $fullText.IndexOf($old) + 1

$target = $tr.Characters($startPos, $old.Length)
$target.Text = $new

$target2 = $tr.Characters($startPos, $new.Length)
$hyperlink = $target2.ActionSettings(1).Hyperlink
$hyperlink.Address = "https://github.com/rfordatascience/rbp_workshop/tree/main/week_05"
